# Breast IG V3 Changes.docx
# Commit: "Changed Breast Composition from hasMember to component"
#
# The target edit adds a brand new bullet/list item to the "Changes"
# list, right after the existing
#   "Changed BiRadsAssessmentCategory from Observation.hasmember to
#    Observation.component in all occurrences."
# item (same ListParagraph style / numId=11 list), with the new text:
#   "Changed Breast Composition Category from hasMember to Component."

$d = $word.ActiveDocument

# Locate the anchor paragraph by its distinctive text, and remember its
# 1-based index within $d.Paragraphs so we can re-fetch a live
# reference after the mutation below (the original COM object can go
# stale once the paragraph collection changes).
$anchorText = "Changed BiRadsAssessmentCategory from Observation.hasmember to Observation.component in all occurrences."
$anchorIndex = 0
$i = 0
foreach ($p in $d.Paragraphs) {
    $i = $i + 1
    if ($p.Range.Text -like "*$anchorText*") {
        $anchorIndex = $i
        break
    }
}

if ($anchorIndex -eq 0) {
    throw "Could not find anchor paragraph"
}

$anchorPara = $d.Paragraphs.Item($anchorIndex)

# Insert a new paragraph right after the anchor paragraph. Word
# automatically carries over the paragraph formatting (style, numPr,
# i.e. ListParagraph / numId=11) from the anchor paragraph.
$anchorPara.Range.InsertParagraphAfter()

# Re-fetch the freshly inserted (currently empty) paragraph by index
# and give it the new text.
$newPara = $d.Paragraphs.Item($anchorIndex + 1)
$newPara.Range.Text = "Changed Breast Composition Category from hasMember to Component."

Write-Host "Inserted new paragraph after anchor (index $anchorIndex)."
